# mappings.xlsx - "cleaned version (manual domain fix)"
#
# 1) Trim stray trailing whitespace off the Domain column values on the
#    SupplierMappings sheet (lukebrands.com, wallisco.com, bylooil.com,
#    world-kinect.com, opisnet.com all had a trailing space).
# 2) Restore the workbook's active sheet/selection to SupplierMappings
#    (cell E11) instead of SupplyMappings (cell A4), which also flips
#    which sheet tab is marked selected.

$wb = $excel.ActiveWorkbook

$supplierSheet = $wb.Worksheets.Item("SupplierMappings")
$supplySheet   = $wb.Worksheets.Item("SupplyMappings")

# --- 1) Clean up the trailing-space domains in column C (rows 2-6) ---
$domainRange = $supplierSheet.Range("C2:C6")
foreach ($cell in $domainRange.Cells) {
    $current = $cell.Value2
    if ($current -ne $null) {
        $cleaned = $current.TrimEnd()
        if ($cleaned -ne $current) {
            $cell.Value = $cleaned
        }
    }
}

# --- 2) Fix up which sheet/cell is active ---
# Make sure SupplyMappings no longer looks like the "current" tab.
$supplySheet.Range("A4").Select()

# SupplierMappings becomes the active/selected tab again, with E11 selected.
$supplierSheet.Activate()
$supplierSheet.Range("E11").Select()
